# Auto-generated script applying scheduled market-data refresh to Atomos_Profits sheets.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H,I,J,K,L,M,N) per leve row
# for the sheets ALC, ARM, BSM, CRP, CUL, GSM, WVR based on refreshed marketboard data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1575.619
$ws.Range("J43").Value = 1711.1177
$ws.Range("L43").Value = 1711.1177
$ws.Range("N43").Value = -1849.1177
$ws.Range("H70").Value = 1642.9546
$ws.Range("I70").Value = 1248.875
$ws.Range("J70").Value = 2693.8333
$ws.Range("K70").Value = 3746.625
$ws.Range("L70").Value = 8081.499899999999
$ws.Range("M70").Value = -3476.625
$ws.Range("N70").Value = -8621.499899999999
$ws.Range("H73").Value = 1642.9546
$ws.Range("I73").Value = 1248.875
$ws.Range("J73").Value = 2693.8333
$ws.Range("K73").Value = 3746.625
$ws.Range("L73").Value = 8081.499899999999
$ws.Range("M73").Value = -2810.625
$ws.Range("N73").Value = -9953.499899999999
$ws.Range("H106").Value = 3123.5
$ws.Range("I106").Value = 1908.75
$ws.Range("J106").Value = 3933.3333
$ws.Range("K106").Value = 1908.75
$ws.Range("L106").Value = 3933.3333
$ws.Range("M106").Value = -1277.75
$ws.Range("N106").Value = -5195.3333
$ws.Range("H107").Value = 640
$ws.Range("I107").Value = 592.8333
$ws.Range("J107").Value = 734.3333
$ws.Range("K107").Value = 592.8333
$ws.Range("L107").Value = 734.3333
$ws.Range("M107").Value = 1327.1667
$ws.Range("N107").Value = -4574.3333
$ws.Range("H116").Value = 4058
$ws.Range("I116").Value = 3400
$ws.Range("J116").Value = 4810
$ws.Range("K116").Value = 3400
$ws.Range("L116").Value = 4810
$ws.Range("M116").Value = 42
$ws.Range("N116").Value = -11694
$ws.Range("H125").Value = 1246.2222
$ws.Range("I125").Value = 1167.3846
$ws.Range("K125").Value = 10506.4614
$ws.Range("M125").Value = -8046.4614
$ws.Range("H129").Value = 1023.1515
$ws.Range("J129").Value = 989.082
$ws.Range("L129").Value = 2967.246
$ws.Range("N129").Value = -12967.246
$ws.Range("H132").Value = 4548057
$ws.Range("I132").Value = 5002389
$ws.Range("J132").Value = 4740
$ws.Range("K132").Value = 15007167
$ws.Range("L132").Value = 14220
$ws.Range("M132").Value = -15004637
$ws.Range("N132").Value = -19280
$ws.Range("H137").Value = 3181.5686
$ws.Range("I137").Value = 3358.8857
$ws.Range("K137").Value = 10076.6571
$ws.Range("M137").Value = -7526.6571
$ws.Range("H141").Value = 410511.72
$ws.Range("I141").Value = 1649.1305
$ws.Range("J141").Value = 1265406.1
$ws.Range("K141").Value = 4947.3915
$ws.Range("L141").Value = 3796218.3
$ws.Range("M141").Value = 232.6085000000003
$ws.Range("N141").Value = -3806578.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2096.6667
$ws.Range("I61").Value = 1279.4546
$ws.Range("J61").Value = 2658.5
$ws.Range("K61").Value = 1279.4546
$ws.Range("L61").Value = 2658.5
$ws.Range("M61").Value = -1067.4546
$ws.Range("N61").Value = -3082.5
$ws.Range("H74").Value = 1186.8948
$ws.Range("I74").Value = 1221.1875
$ws.Range("J74").Value = 1004
$ws.Range("K74").Value = 1221.1875
$ws.Range("L74").Value = 1004
$ws.Range("M74").Value = -347.1875
$ws.Range("N74").Value = -2752
$ws.Range("H77").Value = 1186.8948
$ws.Range("I77").Value = 1221.1875
$ws.Range("J77").Value = 1004
$ws.Range("K77").Value = 6105.9375
$ws.Range("L77").Value = 5020
$ws.Range("M77").Value = -1737.9375
$ws.Range("N77").Value = -13756
$ws.Range("H122").Value = 2099.5957
$ws.Range("I122").Value = 1677.6285
$ws.Range("J122").Value = 3330.3333
$ws.Range("K122").Value = 5032.8855
$ws.Range("L122").Value = 9990.999899999999
$ws.Range("M122").Value = -2582.8855
$ws.Range("N122").Value = -14890.9999
$ws.Range("H132").Value = 1955.5892
$ws.Range("I132").Value = 1611.4186
$ws.Range("J132").Value = 3094
$ws.Range("K132").Value = 4834.2558
$ws.Range("L132").Value = 9282
$ws.Range("M132").Value = -2304.2558
$ws.Range("N132").Value = -14342
$ws.Range("H136").Value = 2096.6667
$ws.Range("I136").Value = 1279.4546
$ws.Range("J136").Value = 2658.5
$ws.Range("K136").Value = 3838.3638
$ws.Range("L136").Value = 7975.5
$ws.Range("M136").Value = -1288.3638
$ws.Range("N136").Value = -13075.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 561.4167
$ws.Range("I80").Value = 618.36365
$ws.Range("J80").Value = 513.2308
$ws.Range("K80").Value = 618.36365
$ws.Range("L80").Value = 513.2308
$ws.Range("M80").Value = 379.63635
$ws.Range("N80").Value = -2509.2308
$ws.Range("H83").Value = 561.4167
$ws.Range("I83").Value = 618.36365
$ws.Range("J83").Value = 513.2308
$ws.Range("K83").Value = 3091.81825
$ws.Range("L83").Value = 2566.154
$ws.Range("M83").Value = 1900.18175
$ws.Range("N83").Value = -12550.154
$ws.Range("H134").Value = 2038.375
$ws.Range("I134").Value = 1442.138
$ws.Range("J134").Value = 3610.2727
$ws.Range("K134").Value = 4326.414
$ws.Range("L134").Value = 10830.8181
$ws.Range("M134").Value = -1791.414
$ws.Range("N134").Value = -15900.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3357.2092
$ws.Range("I31").Value = 2595.7083
$ws.Range("J31").Value = 4319.1055
$ws.Range("K31").Value = 2595.7083
$ws.Range("L31").Value = 4319.1055
$ws.Range("M31").Value = -2300.7083
$ws.Range("N31").Value = -4909.1055
$ws.Range("H34").Value = 3357.2092
$ws.Range("I34").Value = 2595.7083
$ws.Range("J34").Value = 4319.1055
$ws.Range("K34").Value = 2595.7083
$ws.Range("L34").Value = 4319.1055
$ws.Range("M34").Value = -2393.7083
$ws.Range("N34").Value = -4723.1055
$ws.Range("H99").Value = 1999.0769
$ws.Range("I99").Value = 1284.8572
$ws.Range("J99").Value = 2832.3333
$ws.Range("K99").Value = 1284.8572
$ws.Range("L99").Value = 2832.3333
$ws.Range("M99").Value = 213.1428000000001
$ws.Range("N99").Value = -5828.3333
$ws.Range("H107").Value = 1280.2727
$ws.Range("I107").Value = 421.5
$ws.Range("K107").Value = 421.5
$ws.Range("M107").Value = 1498.5
$ws.Range("H126").Value = 1999.0769
$ws.Range("I126").Value = 1284.8572
$ws.Range("J126").Value = 2832.3333
$ws.Range("K126").Value = 3854.5716
$ws.Range("L126").Value = 8496.999899999999
$ws.Range("M126").Value = -1384.5716
$ws.Range("N126").Value = -13436.9999
$ws.Range("H132").Value = 3141.32
$ws.Range("I132").Value = 2568.9333
$ws.Range("J132").Value = 3999.9
$ws.Range("K132").Value = 7706.7999
$ws.Range("L132").Value = 11999.7
$ws.Range("M132").Value = -5176.7999
$ws.Range("N132").Value = -17059.7
$ws.Range("H134").Value = 2468.8096
$ws.Range("I134").Value = 828.75
$ws.Range("J134").Value = 4655.5557
$ws.Range("K134").Value = 2486.25
$ws.Range("L134").Value = 13966.6671
$ws.Range("M134").Value = 48.75
$ws.Range("N134").Value = -19036.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 902
$ws.Range("J113").Value = 950.6316
$ws.Range("L113").Value = 2851.8948
$ws.Range("N113").Value = -7191.8948
$ws.Range("H131").Value = 1052.9387
$ws.Range("I131").Value = 5050
$ws.Range("J131").Value = 969.6667
$ws.Range("K131").Value = 15150
$ws.Range("L131").Value = 2909.0001
$ws.Range("M131").Value = -10110
$ws.Range("N131").Value = -12989.0001
$ws.Range("H140").Value = 1614.5862
$ws.Range("I140").Value = 828.4091
$ws.Range("J140").Value = 4085.4285
$ws.Range("K140").Value = 2485.2273
$ws.Range("L140").Value = 12256.2855
$ws.Range("M140").Value = 2694.7727
$ws.Range("N140").Value = -22616.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 26672.65
$ws.Range("I102").Value = 2760.524
$ws.Range("J102").Value = 49497.863
$ws.Range("K102").Value = 2760.524
$ws.Range("L102").Value = 49497.863
$ws.Range("M102").Value = -1138.524
$ws.Range("N102").Value = -52741.863
$ws.Range("H132").Value = 4617.5835
$ws.Range("I132").Value = 8137.3335
$ws.Range("J132").Value = 3444.3333
$ws.Range("K132").Value = 24412.0005
$ws.Range("L132").Value = 10332.9999
$ws.Range("M132").Value = -21882.0005
$ws.Range("N132").Value = -15392.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 70004
$ws.Range("J8").Value = 70004
$ws.Range("L8").Value = 70004
$ws.Range("N8").Value = -70284
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H126").Value = 2942988.5
$ws.Range("I126").Value = 1414.2727
$ws.Range("J126").Value = 8335875
$ws.Range("K126").Value = 4242.8181
$ws.Range("L126").Value = 25007625
$ws.Range("M126").Value = -1772.8181
$ws.Range("N126").Value = -25012565
$ws.Range("H138").Value = 29459.428
$ws.Range("J138").Value = 29459.428
$ws.Range("L138").Value = 29459.428
$ws.Range("N138").Value = -39739.428
